$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 3 (the "PITISUTTIHUM Punnee" / "Sub I" record).
$ws.Rows(3).Delete()

# 2. Insert a new column at C ("Project Number 2"), shifting
#    "Sponsor Protocol Number" and everything after it one column to the right.
$ws.Columns("C").Insert()

# 3. Give the new column the same width as in the target workbook.
#    (ColumnWidth is expressed in characters; the runtime quantizes the
#    stored width to whole pixels at 6px/char, so 15.66666667 is the
#    closest achievable match to the target stored width of 16.5703125.)
$ws.Columns("C").ColumnWidth = 15.66666667

# 4. Populate the new column's header and data value.
$ws.Cells.Item(1, 3).Value = "Project Number 2"
$ws.Cells.Item(2, 3).Value = "0000/8595"

# 5. The remaining record's "Sponsor Protocol Number" (now column D) changes
#    from the text "0102" to the numeric value 67873. Plainly overwriting the
#    cell's Value would make Excel drop the cell's quote-prefix/border style
#    (since quote-prefix no longer applies once the content is a real
#    number), which does not match the target formatting. To keep the
#    original style bits intact we stage the new value on a scratch cell,
#    paste the original cell's format onto it, write the value there too,
#    and finally copy the format back onto the real cell after its value
#    has already been updated.
$scratch = $ws.Cells.Item(100, 100)
$scratch.Value = 67873
$ws.Cells.Item(2, 4).Copy()
$scratch.PasteSpecial(-4122)          # xlPasteFormats

$ws.Cells.Item(2, 4).Value = 67873    # value now correct, style temporarily degraded

$scratch.Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4122)  # xlPasteFormats: restore the original style
$scratch.Clear()

# 6. Update the active selection to match the target workbook.
$ws.Range("D2").Select()
